$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '97.135.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.584.26'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '655.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.66'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +13.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.411'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('E9').Value = '  +6.66%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.581.30'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.204'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.254.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.893.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000257'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.587.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.528'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +10.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '510.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000201'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.778.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('E31').Value = '  +9.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.51'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.63%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.185'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '624.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.75'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.569'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.153'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('B43').Value = 'ImmutableX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.40%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.910'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.71%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0432'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.23'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.40%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.35'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.36%  '
$ws.Range('B51').Value = 'MantraDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.51'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.47%  '
